# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型"
# sheets to the freshly generated numbers.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---------------------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value  = 1584
$ws1.Range("F8").Value  = 141
$ws1.Range("F10").Value = 8420
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 16
$ws1.Range("F15").Value = 1301
$ws1.Range("F16").Value = 66
$ws1.Range("F18").Value = 28
$ws1.Range("F19").Value = 9043
$ws1.Range("F20").Value = 147
$ws1.Range("F21").Value = 88
$ws1.Range("F22").Value = 204
$ws1.Range("F25").Value = 5871
$ws1.Range("F26").Value = 1032
$ws1.Range("F27").Value = 43
$ws1.Range("F29").Value = 85

# --- Sheet: 全部类型 ------------------------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 1584
$ws4.Range("F9").Value  = 141
$ws4.Range("F11").Value = 8420
$ws4.Range("F13").Value = 57
$ws4.Range("F14").Value = 16
$ws4.Range("F16").Value = 1301
$ws4.Range("F17").Value = 66
$ws4.Range("F19").Value = 28
$ws4.Range("F22").Value = 9043
$ws4.Range("F23").Value = 147
$ws4.Range("F24").Value = 88
$ws4.Range("F25").Value = 204
$ws4.Range("F28").Value = 5871
$ws4.Range("F29").Value = 1032
$ws4.Range("F30").Value = 43
$ws4.Range("F32").Value = 85
